# Loan RBI, Variable Instalments
# Insert a new (currently blank) column before the "Late" column on the
# "Repayment schedule" sheet, and make that sheet the active one with
# cell K15 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a blank column at N (pushes Late/Date/Outstanding one column right)
$ws.Columns("N").Insert()

# Match the width of the neighbouring "Outstanding"/In Advance column so the
# new blank column doesn't look squeezed.
$ws.Columns("N").ColumnWidth = 9.83

# Make "Repayment schedule" the active sheet and select K15, as in the
# edited workbook.
$ws.Activate()
$ws.Range("K15").Select()
